$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Update URL value (row 2, column B)
$ws1.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-diseasestage-codes"

# Update Experimental value (row 7, column B) - set to "true" as text
$ws1.Range("B7").Value = "'true"

# Update Date value (row 8, column B)
$ws1.Range("B8").Value = "2024-02-19T18:37:26-06:00"
